# Update the "2019 K-Means Clustering Data" table on the slide with the
# cluster-statistics table (Slide 26, Shape 2 - a graphicFrame containing
# a table with 3 rows x 11 columns).
#
# Columns: 1 Cluster, 2 Teams, 3 MeanX2014HF, 4 MeanX2015HF, 5 MeanX2016HF,
#          6 MeanX2017HF, 7 MeanX2018HF, 8 MeanX2019HF, 9 MeanPopulation,
#          10 MeanSalary, 11 MeanESPNRating

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(26)
$shape = $s.Shapes.Item(2)
$tbl = $shape.Table

# Row 2 (Cluster 1) updates
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "36"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "1.9"
$tbl.Cell(2, 4).Shape.TextFrame.TextRange.Text = "2.3"
$tbl.Cell(2, 5).Shape.TextFrame.TextRange.Text = "2.3"
$tbl.Cell(2, 6).Shape.TextFrame.TextRange.Text = "1.8"
$tbl.Cell(2, 7).Shape.TextFrame.TextRange.Text = "1.5"
$tbl.Cell(2, 8).Shape.TextFrame.TextRange.Text = "1.1"
$tbl.Cell(2, 9).Shape.TextFrame.TextRange.Text = "1.8"
$tbl.Cell(2, 10).Shape.TextFrame.TextRange.Text = "14.0"
$tbl.Cell(2, 11).Shape.TextFrame.TextRange.Text = "3.8"

# Row 3 (Cluster 2) updates
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "87"
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "0.5"
$tbl.Cell(3, 4).Shape.TextFrame.TextRange.Text = "0.3"
$tbl.Cell(3, 6).Shape.TextFrame.TextRange.Text = "0.6"
$tbl.Cell(3, 7).Shape.TextFrame.TextRange.Text = "0.7"
$tbl.Cell(3, 8).Shape.TextFrame.TextRange.Text = "0.9"
$tbl.Cell(3, 10).Shape.TextFrame.TextRange.Text = "16.7"
